$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file path strings in column A (rows 2 and 3)
$ws.Range("A2").Value = "/home/daniel/Spike Data/Matlab files/Exp 19 baseline data.mat"
$ws.Range("A3").Value = "/home/daniel/Spike Data/Matlab files/Exp 27 unit 1 data.mat"

# Move the active cell selection to A6
$ws.Range("A6").Select()

# Adjust the tab ratio (split between sheet tabs and horizontal scroll bar).
# Excel's COM TabRatio is a 0..1 fraction of the OOXML bookViews@tabRatio
# permille value (e.g. tabRatio="76" <-> TabRatio = 0.076).
$excel.ActiveWindow.TabRatio = 0.076
